# Update Result/Date columns (A2/B2) on each sheet, reflecting a newer
# Katalon test run (ABP Object repository, Keyword Classes and ABP test
# cases added -> re-run timestamps/results recorded in the bootstrap data).

$wb = $excel.ActiveWorkbook

$updates = @{
    "PayNowNoCFPC_27"   = @{ Result = "Pass"; Date = "Thu Jun 05 11:59:09 IST 2025" }
    "PayNowNoCFPS_27"   = @{ Result = "Pass"; Date = "Thu Jun 05 11:59:44 IST 2025" }
    "PayNowNoCFCorp_27" = @{ Result = "Pass"; Date = "Thu Jun 05 11:58:14 IST 2025" }
    "PayNowSCFPC_27"    = @{ Result = "Pass"; Date = "Thu Jun 05 12:02:00 IST 2025" }
    "PayNowSCFPS_27"    = @{ Result = "Pass"; Date = "Thu Jun 05 12:02:55 IST 2025" }
    "PayNowSCFCorp_27"  = @{ Result = "Pass"; Date = "Thu Jun 05 12:01:26 IST 2025" }
    "PayNowDCFPC_27"    = @{ Result = "Pass"; Date = "Thu Jun 05 11:55:17 IST 2025" }
    "PayNowDCFPS_27"    = @{ Result = "Pass"; Date = "Thu Jun 05 11:56:17 IST 2025" }
    "PayNowDCFCorp_27"  = @{ Result = "Pass"; Date = "Tue Jun 17 11:13:14 IST 2025" }
    "CCDeferredPS_27"   = @{ Result = "Fail"; Date = "Thu Jun 05 11:43:18 IST 2025" }
    "CCDeferredPC_27"   = @{ Result = "Fail"; Date = "Tue Jun 10 14:25:57 IST 2025" }
    "CCDeferredCorp_27" = @{ Result = "Pass"; Date = "Tue Jun 17 11:38:24 IST 2025" }
    "CMCAutopayPC_27"   = @{ Result = "Pass"; Date = "Thu Jun 05 11:49:29 IST 2025" }
    "CMCAutopayCorp_27" = @{ Result = "Pass"; Date = "Thu Jun 05 11:47:38 IST 2025" }
    "CMCAutopayPS_27"   = @{ Result = "Pass"; Date = "Thu Jun 05 11:51:23 IST 2025" }
}

foreach ($ws in $wb.Worksheets) {
    $info = $updates[$ws.Name]
    if ($info -ne $null) {
        $ws.Range("A2").Value = $info.Result
        $ws.Range("B2").Value = $info.Date
    }
}
